# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet (zh-cn/de-de status cells) and on each language sheet.
# - The de-de/zh-cn "Latest Handback DateTime" timestamps advance to reflect
#   the new handback pass.
# - The stale-handback-version warning in "Error Detail" is cleared now that
#   the handback is in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
# The status columns widen to fit the longer "Handed back" text.
$overview.Columns.Item(5).ColumnWidth = 29.1666666
$overview.Columns.Item(6).ColumnWidth = 29.1666666

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-24 16:50:53"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.1666666
$zhcn.Columns.Item(16).ColumnWidth = 12.8333333

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-24 16:51:03"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.1666666
$dede.Columns.Item(16).ColumnWidth = 12.8333333
